$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.914.61"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "2.930.23"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "374.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.40%  "
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("D13").Value = "3.394.92"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.83%  "
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").Value = "2.929.49"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").Value = "50.897.93"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("E19").Value = "  -6.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.92%  "
$ws.Range("E21").Value = "  -3.15%  "
$ws.Range("D22").Value = "0.0₃0952"
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "263.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("B26").Value = "Filecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.74%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.63%  "
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0447"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  -3.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.51%  "
$ws.Range("E39").Value = "  -2.49%  "
$ws.Range("E40").Value = "  -1.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.77%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.279"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "1.965.23"
$ws.Range("E49").Value = "  -3.40%  "
$ws.Range("E50").Value = "  -3.17%  "
$ws.Range("E51").Value = "  -2.59%  "
